$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 847, pushing existing rows 847-888 down to 848-889.
$ws.Rows.Item(847).Insert()

# Populate the newly inserted row 847 with its data.
# Column A holds a date stored as plain text (e.g. "2026/02/24"), so force
# a text number format first to prevent Excel from auto-converting the
# string into a real date value.
$ws.Range("A847").NumberFormat = "@"
$ws.Range("A847").Value = "2026/02/24"

$ws.Range("B847").Value = "火"
$ws.Range("C847").Value = 12
$ws.Range("D847").Value = 105
